$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1361.6923
$ws.Range("J19").Value = 1323.3334
$ws.Range("L19").Value = 1323.3334
$ws.Range("N19").Value = -1673.3334
$ws.Range("H32").Value = 4627.3706
$ws.Range("J32").Value = 2086.3157
$ws.Range("L32").Value = 2086.3157
$ws.Range("N32").Value = -2738.3157
$ws.Range("H42").Value = 162.125
$ws.Range("I42").Value = 42.42857
$ws.Range("J42").Value = 1000
$ws.Range("K42").Value = 127.28571
$ws.Range("L42").Value = 3000
$ws.Range("M42").Value = 102.71429
$ws.Range("N42").Value = -3460
$ws.Range("H106").Value = 1859
$ws.Range("I106").Value = 1859
$ws.Range("K106").Value = 1859
$ws.Range("M106").Value = -1228
$ws.Range("H116").Value = 17416.8
$ws.Range("J116").Value = 11856.25
$ws.Range("L116").Value = 11856.25
$ws.Range("N116").Value = -18740.25
$ws.Range("H135").Value = 3665.9375
$ws.Range("I135").Value = 2877.7273
$ws.Range("J135").Value = 5400
$ws.Range("K135").Value = 25899.5457
$ws.Range("L135").Value = 48600
$ws.Range("M135").Value = -23364.5457
$ws.Range("N135").Value = -53670
$ws.Range("H137").Value = 27033
$ws.Range("J137").Value = 3559
$ws.Range("L137").Value = 10677
$ws.Range("N137").Value = -15777
$ws.Range("H138").Value = 34756.195
$ws.Range("J138").Value = 74408.36
$ws.Range("L138").Value = 223225.08
$ws.Range("N138").Value = -233505.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29497.838
$ws.Range("I32").Value = 33950.78
$ws.Range("K32").Value = 33950.78
$ws.Range("M32").Value = -33663.78
$ws.Range("H61").Value = 13828.111
$ws.Range("I61").Value = 3300
$ws.Range("J61").Value = 15144.125
$ws.Range("K61").Value = 3300
$ws.Range("L61").Value = 15144.125
$ws.Range("M61").Value = -3088
$ws.Range("N61").Value = -15568.125
$ws.Range("H97").Value = 973.1724
$ws.Range("I97").Value = 783.12
$ws.Range("J97").Value = 2161
$ws.Range("K97").Value = 783.12
$ws.Range("L97").Value = 2161
$ws.Range("M97").Value = -287.12
$ws.Range("N97").Value = -3153
$ws.Range("H110").Value = 391.9
$ws.Range("I110").Value = 391.9
$ws.Range("K110").Value = 391.9
$ws.Range("M110").Value = 1653.1
$ws.Range("H122").Value = 2683.195
$ws.Range("J122").Value = 2901.25
$ws.Range("L122").Value = 8703.75
$ws.Range("N122").Value = -13603.75
$ws.Range("H136").Value = 13828.111
$ws.Range("I136").Value = 3300
$ws.Range("J136").Value = 15144.125
$ws.Range("K136").Value = 9900
$ws.Range("L136").Value = 45432.375
$ws.Range("M136").Value = -7350
$ws.Range("N136").Value = -50532.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H64").Value = 13720.5
$ws.Range("I64").Value = 735.5
$ws.Range("K64").Value = 735.5
$ws.Range("M64").Value = -510.5
$ws.Range("H67").Value = 13720.5
$ws.Range("I67").Value = 735.5
$ws.Range("K67").Value = 735.5
$ws.Range("M67").Value = 44.5
$ws.Range("H86").Value = 1828.9
$ws.Range("I86").Value = 1631.5
$ws.Range("K86").Value = 1631.5
$ws.Range("M86").Value = -508.5
$ws.Range("H89").Value = 1828.9
$ws.Range("I89").Value = 1631.5
$ws.Range("K89").Value = 8157.5
$ws.Range("M89").Value = -2541.5
$ws.Range("H105").Value = 1357.973
$ws.Range("I105").Value = 1030.2727
$ws.Range("J105").Value = 1838.6
$ws.Range("K105").Value = 1030.2727
$ws.Range("L105").Value = 1838.6
$ws.Range("M105").Value = 716.7273
$ws.Range("N105").Value = -5332.6
$ws.Range("H107").Value = 2879.0833
$ws.Range("I107").Value = 2745.5386
$ws.Range("K107").Value = 2745.5386
$ws.Range("M107").Value = -825.5385999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3703.182
$ws.Range("I105").Value = 1570
$ws.Range("J105").Value = 4922.143
$ws.Range("K105").Value = 1570
$ws.Range("L105").Value = 4922.143
$ws.Range("M105").Value = 177
$ws.Range("N105").Value = -8416.143
$ws.Range("H107").Value = 817.94446
$ws.Range("I107").Value = 747.5454999999999
$ws.Range("J107").Value = 928.5714
$ws.Range("K107").Value = 747.5454999999999
$ws.Range("L107").Value = 928.5714
$ws.Range("M107").Value = 1172.4545
$ws.Range("N107").Value = -4768.5714

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1040.9
$ws.Range("J98").Value = 1117.6666
$ws.Range("L98").Value = 3352.9998
$ws.Range("N98").Value = -6348.9998
$ws.Range("H107").Value = 2182.611
$ws.Range("I107").Value = 3934.25
$ws.Range("J107").Value = 781.3
$ws.Range("K107").Value = 11802.75
$ws.Range("L107").Value = 2343.9
$ws.Range("M107").Value = -9882.75
$ws.Range("N107").Value = -6183.9
$ws.Range("H123").Value = 2629
$ws.Range("I123").Value = 2405
$ws.Range("J123").Value = 3525
$ws.Range("K123").Value = 7215
$ws.Range("L123").Value = 10575
$ws.Range("M123").Value = -4765
$ws.Range("N123").Value = -15475
$ws.Range("H140").Value = 2869.875
$ws.Range("I140").Value = 2869.875
$ws.Range("K140").Value = 8609.625
$ws.Range("M140").Value = -3429.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 496.18518
$ws.Range("I2").Value = 628
$ws.Range("K2").Value = 628
$ws.Range("M2").Value = -515
$ws.Range("H96").Value = 49500
$ws.Range("J96").Value = 49500
$ws.Range("L96").Value = 49500
$ws.Range("N96").Value = -54992
$ws.Range("H126").Value = 2638.5
$ws.Range("I126").Value = 2158.2856
$ws.Range("K126").Value = 6474.8568
$ws.Range("M126").Value = -4004.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 125000
$ws.Range("J63").Value = 125000
$ws.Range("L63").Value = 125000
$ws.Range("N63").Value = -126498
$ws.Range("H66").Value = 125000
$ws.Range("J66").Value = 125000
$ws.Range("L66").Value = 375000
$ws.Range("N66").Value = -382488
$ws.Range("H81").Value = 150000
$ws.Range("J81").Value = 150000
$ws.Range("L81").Value = 150000
$ws.Range("N81").Value = -151996
$ws.Range("H84").Value = 150000
$ws.Range("J84").Value = 150000
$ws.Range("L84").Value = 450000
$ws.Range("N84").Value = -459984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 19875.166
$ws.Range("J45").Value = 19875.166
$ws.Range("L45").Value = 19875.166
$ws.Range("N45").Value = -20857.166
$ws.Range("H62").Value = 4066.3333
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 4066.3333
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240
$ws.Range("H107").Value = 1011.7059
$ws.Range("I107").Value = 965.4545000000001
$ws.Range("K107").Value = 2896.3635
$ws.Range("M107").Value = -976.3635000000004
$ws.Range("H122").Value = 62677.582
$ws.Range("I122").Value = 74398.45
$ws.Range("J122").Value = 4073.25
$ws.Range("K122").Value = 223195.35
$ws.Range("L122").Value = 12219.75
$ws.Range("M122").Value = -220745.35
$ws.Range("N122").Value = -17119.75
$ws.Range("H135").Value = 63499.5
$ws.Range("J135").Value = 63499.5
$ws.Range("L135").Value = 63499.5
$ws.Range("N135").Value = -73639.5
$ws.Range("H136").Value = 19766.156
$ws.Range("I136").Value = 24851.209
$ws.Range("K136").Value = 74553.62699999999
$ws.Range("M136").Value = -72003.62699999999
